# Edit Config.xlsx: leave the bot ready to add mailbox (O365) credentials.
# Commit message: "Dejo listo el bot para poder agregar credenciales de buzon"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# --- Update existing parameter values -------------------------------------

# CarpetaCompartida: point to the real shared network folder instead of the
# old local test path.
$ws.Cells.Item(6, 2).Value = "//10.238.99.5/temporales/ROCKETBOT/Compartida"

# RutaArchivoBase: drop the "Copia de " prefix from the input file name.
$ws.Cells.Item(7, 2).Value = "//10.238.99.5/temporales/ROCKETBOT/Inusmo busqueda masiva - UT Ticenergi.xlsx"

# --- Add new parameter rows for the O365 mailbox credential ---------------

# Row 10: CredencialO365
$ws.Cells.Item(10, 1).Value = "CredencialO365"
$ws.Cells.Item(10, 2).Value = "O365"
$ws.Cells.Item(10, 3).Value = "Valor"
$ws.Cells.Item(10, 4).Value = "Credenciales O365"
$ws.Rows.Item(10).RowHeight = 15

# Row 11: TenantId (Azure AD tenant identifier)
$ws.Cells.Item(11, 1).Value = "TenantId"
$ws.Cells.Item(11, 2).NumberFormat = "0.00E+00"
$ws.Cells.Item(11, 2).Value = "024e2966-f777-426e-9ffb-80231972a6b1"
$ws.Cells.Item(11, 3).Value = "Valor"
$ws.Cells.Item(11, 4).Value = "Tenant ID de Azure"
$ws.Rows.Item(11).RowHeight = 15

# Reflect where the editor left the selection/scroll after the edits.
$ws.Range("B13").Select()
